$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.901.82"
$ws.Range("E2").Value = "  -0.44%  "
$ws.Range("D3").Value = "2.789.34"
$ws.Range("E3").Value = "  -1.90%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'360.63"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("D6").Value = "'109.67"
$ws.Range("E6").Value = "  -3.68%  "
$ws.Range("D7").Value = "'0.558"
$ws.Range("E7").Value = "  -2.16%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'0.590"
$ws.Range("D10").Value = "'39.95"
$ws.Range("E10").Value = "  -4.23%  "
$ws.Range("D11").Value = "'0.0845"
$ws.Range("E11").Value = "  -2.12%  "
$ws.Range("E12").Value = "  +1.06%  "
$ws.Range("D13").Value = "'19.48"
$ws.Range("E13").Value = "  -2.64%  "
$ws.Range("D14").Value = "'7.55"
$ws.Range("E14").Value = "  -2.87%  "
$ws.Range("D15").Value = "3.224.86"
$ws.Range("E15").Value = "  -1.92%  "
$ws.Range("D16").Value = "2.787.36"
$ws.Range("E16").Value = "  -1.72%  "
$ws.Range("D17").Value = "'0.938"
$ws.Range("E17").Value = "  +2.54%  "
$ws.Range("D18").Value = "51.878.82"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").Value = "'7.49"
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("D20").Value = "'3.10"
$ws.Range("E20").Value = "  -2.41%  "
$ws.Range("D21").Value = "'13.13"
$ws.Range("E21").Value = "  -3.45%  "
$ws.Range("E22").Value = "  -1.98%  "
$ws.Range("D23").Value = "'70.28"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").Value = "'269.84"
$ws.Range("E24").Value = "  +1.00%  "
$ws.Range("D25").Value = "'2.75"
$ws.Range("E25").Value = "  -2.90%  "
$ws.Range("E26").Value = "  -2.61%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("D28").Value = "'0.161"
$ws.Range("E28").Value = "  +14.50%  "
$ws.Range("D29").Value = "'10.27"
$ws.Range("E29").Value = "  -1.62%  "
$ws.Range("E30").Value = "  +1.49%  "
$ws.Range("D31").Value = "'0.0471"
$ws.Range("E31").Value = "  +2.74%  "
$ws.Range("D32").Value = "'51.93"
$ws.Range("E32").Value = "  -2.25%  "
$ws.Range("D33").Value = "'34.08"
$ws.Range("E33").Value = "  -0.23%  "
$ws.Range("E34").Value = "  -3.02%  "
$ws.Range("D35").Value = "'0.0842"
$ws.Range("E35").Value = "  -0.15%  "
$ws.Range("D36").Value = "'5.23"
$ws.Range("E36").Value = "  -2.87%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D38").Value = "'18.96"
$ws.Range("E38").Value = "  +3.41%  "
$ws.Range("E39").Value = "  -2.70%  "
$ws.Range("E40").Value = "  -4.30%  "
$ws.Range("D41").Value = "'2.59"
$ws.Range("E41").Value = "  +1.27%  "
$ws.Range("E42").Value = "  -1.90%  "
$ws.Range("D43").Value = "'2.24"
$ws.Range("E43").Value = "  -0.86%  "
$ws.Range("D44").Value = "'119.73"
$ws.Range("E44").Value = "  -6.53%  "
$ws.Range("D45").Value = "'21.87"
$ws.Range("E45").Value = "  -8.87%  "
$ws.Range("D46").Value = "2.083.07"
$ws.Range("E46").Value = "  -1.90%  "
$ws.Range("E47").Value = "  -4.35%  "
$ws.Range("E48").Value = "  -1.86%  "
$ws.Range("D49").Value = "'5.80"
$ws.Range("E49").Value = "  -0.89%  "
$ws.Range("E50").Value = "  -4.79%  "
$ws.Range("D51").Value = "'8.86"
$ws.Range("E51").Value = "  -2.00%  "
